$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New bold "reviewed" date entry (row 9) ---
$ws.Range("B6").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = 42634
$ws.Range("A9").Font.Bold = $true

# --- New deliverable rows (Hardware Block Diagram / Product Architecture / Major Components BOM) ---
$ws.Range("A10").Value = "Hardware Block Diagram"
$ws.Range("B6").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = 42639

$ws.Range("A11").Value = "Product Architecture"
$ws.Range("B6").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = 42639

$ws.Range("A12").Value = "Major Components BOM"
$ws.Range("B6").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value = 42639

# --- New "Notes" column header on the existing header row (row 5), bold to match "Assigned"/"Due" ---
$ws.Range("C5").Value = "Notes"
$ws.Range("C5").Font.Bold = $true

# Make the Assigned-date cell (A5) bold as well, matching the reviewed layout.
$ws.Range("A5").Font.Bold = $true

# --- Review note for the first deliverable ---
$ws.Range("C12").Value = "I will order parts on Monday to make sure they are in house by Saturday. If you do not get me the parts order by then you may not be able to build you project. "

# --- Column A got a bit wider to fit the new content ---
$ws.Columns("A").ColumnWidth = 21.5
